{"js": "// Map of old cell text -> new cell text, derived from the authoritative diff.\n// Keys are unique across the document, so a single lookup pass (using the\n// text captured BEFORE any edits are applied) is safe even though one of\n// the new values (\"34\u00f77=4, 6\") happens to equal another cell's old value.\nconst replacements = {\n  \"88\u00f78=11, 0\": \"28\u00f75=5, 3\",\n  \"80\u00f78=10, 0\": \"77\u00f79=8, 5\",\n  \"61\u00f72=30, 1\": \"62\u00f79=6, 8\",\n  \"51\u00f72=25, 1\": \"92\u00f72=46, 0\",\n  \"57\u00f74=14, 1\": \"20\u00f79=2, 2\",\n  \"62\u00f76=10, 2\": \"34\u00f77=4, 6\",\n  \"50\u00f72=25, 0\": \"41\u00f75=8, 1\",\n  \"58\u00f76=9, 4\": \"53\u00f73=17, 2\",\n  \"68\u00f78=8, 4\": \"99\u00f77=14, 1\",\n  \"90\u00f78=11, 2\": \"51\u00f78=6, 3\",\n  \"17\u00f73=5, 2\": \"59\u00f77=8, 3\",\n  \"40\u00f76=6, 4\": \"39\u00f72=19, 1\",\n  \"70\u00f79=7, 7\": \"45\u00f74=11, 1\",\n  \"97\u00f78=12, 1\": \"29\u00f75=5, 4\",\n  \"84\u00f75=16, 4\": \"24\u00f78=3, 0\",\n  \"79\u00f79=8, 7\": \"67\u00f73=22, 1\",\n  \"53\u00f78=6, 5\": \"93\u00f75=18, 3\",\n  \"72\u00f79=8, 0\": \"71\u00f75=14, 1\",\n  \"87\u00f77=12, 3\": \"26\u00f77=3, 5\",\n  \"58\u00f79=6, 4\": \"41\u00f79=4, 5\",\n  \"53\u00f74=13, 1\": \"85\u00f78=10, 5\",\n  \"69\u00f78=8, 5\": \"53\u00f79=5, 8\",\n  \"51\u00f77=7, 2\": \"69\u00f77=9, 6\",\n  \"34\u00f77=4, 6\": \"36\u00f75=7, 1\",\n  \"52\u00f76=8, 4\": \"48\u00f72=24, 0\",\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst columnCount = table.values.length > 0 ? table.values[0].length : 0;\n\n// Collect every paragraph in every cell, and load its text in one batch.\nconst paragraphs = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items/text\");\n    paragraphs.push(cellParagraphs);\n  }\n}\nawait context.sync();\n\n// Now that every paragraph's original text has been captured, apply the\n// replacements. Doing the lookup against the pre-loaded text (rather than\n// re-reading after each mutation) avoids any ordering hazard even though\n// one replacement's new value equals another cell's original value.\nfor (const cellParagraphs of paragraphs) {\n  for (const p of cellParagraphs.items) {\n    const newText = replacements[p.text];\n    if (newText !== undefined) {\n      p.insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Map of old cell text -> new cell text, derived from the authoritative diff.\n# Keys are unique across the document, so a single lookup pass (against text\n# captured BEFORE any edits are applied) is safe even though one of the new\n# values (\"34\u00f77=4, 6\") happens to equal another cell's original value.\n$replacements = @{\n    \"88\u00f78=11, 0\" = \"28\u00f75=5, 3\"\n    \"80\u00f78=10, 0\" = \"77\u00f79=8, 5\"\n    \"61\u00f72=30, 1\" = \"62\u00f79=6, 8\"\n    \"51\u00f72=25, 1\" = \"92\u00f72=46, 0\"\n    \"57\u00f74=14, 1\" = \"20\u00f79=2, 2\"\n    \"62\u00f76=10, 2\" = \"34\u00f77=4, 6\"\n    \"50\u00f72=25, 0\" = \"41\u00f75=8, 1\"\n    \"58\u00f76=9, 4\" = \"53\u00f73=17, 2\"\n    \"68\u00f78=8, 4\" = \"99\u00f77=14, 1\"\n    \"90\u00f78=11, 2\" = \"51\u00f78=6, 3\"\n    \"17\u00f73=5, 2\" = \"59\u00f77=8, 3\"\n    \"40\u00f76=6, 4\" = \"39\u00f72=19, 1\"\n    \"70\u00f79=7, 7\" = \"45\u00f74=11, 1\"\n    \"97\u00f78=12, 1\" = \"29\u00f75=5, 4\"\n    \"84\u00f75=16, 4\" = \"24\u00f78=3, 0\"\n    \"79\u00f79=8, 7\" = \"67\u00f73=22, 1\"\n    \"53\u00f78=6, 5\" = \"93\u00f75=18, 3\"\n    \"72\u00f79=8, 0\" = \"71\u00f75=14, 1\"\n    \"87\u00f77=12, 3\" = \"26\u00f77=3, 5\"\n    \"58\u00f79=6, 4\" = \"41\u00f79=4, 5\"\n    \"53\u00f74=13, 1\" = \"85\u00f78=10, 5\"\n    \"69\u00f78=8, 5\" = \"53\u00f79=5, 8\"\n    \"51\u00f77=7, 2\" = \"69\u00f77=9, 6\"\n    \"34\u00f77=4, 6\" = \"36\u00f75=7, 1\"\n    \"52\u00f76=8, 4\" = \"48\u00f72=24, 0\"\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# First pass: capture every cell's current text (trimming the trailing\n# cell-mark / end-of-cell characters Word appends to Range.Text), so the\n# replacement lookups below all use the document's ORIGINAL content.\n$cellsToUpdate = @()\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $raw = $cell.Range.Text\n        $clean = $raw.TrimEnd([char]7).TrimEnd([char]13)\n        if ($replacements.ContainsKey($clean)) {\n            $cellsToUpdate += [PSCustomObject]@{ Row = $r; Col = $c; NewText = $replacements[$clean] }\n        }\n    }\n}\n\n# Second pass: apply the replacements. Assigning Range.Text keeps the\n# existing run/paragraph formatting (font, size, alignment) intact.\nforeach ($item in $cellsToUpdate) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $cell.Range.Text = $item.NewText\n}\n"}
